$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D are written with a leading apostrophe so Excel
# stores them as text (quote-prefixed), matching the source data which
# contains values such as thousands-dotted numbers that must remain literal
# strings rather than being auto-converted to numbers.

# Row 2
$ws.Range("D2").Value = "'74.160.37"
$ws.Range("E2").Value = "  +8.97%  "

# Row 3
$ws.Range("D3").Value = "'2.592.84"
$ws.Range("E3").Value = "  +7.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'585.42"
$ws.Range("E5").Value = "  +5.50%  "

# Row 6
$ws.Range("D6").Value = "'181.68"
$ws.Range("E6").Value = "  +14.33%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.208"
$ws.Range("E8").Value = "  +27.57%  "

# Row 9
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +5.24%  "

# Row 10
$ws.Range("D10").Value = "'2.584.91"
$ws.Range("E10").Value = "  +7.21%  "

# Row 11
$ws.Range("E11").Value = "  -0.60%  "

# Row 12
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +7.75%  "

# Row 13
$ws.Range("E13").Value = "  +3.36%  "

# Row 14
$ws.Range("D14").Value = "'0.0000194"
$ws.Range("E14").Value = "  +11.62%  "

# Row 15
$ws.Range("D15").Value = "'74.479.40"
$ws.Range("E15").Value = "  +9.61%  "

# Row 16
$ws.Range("D16").Value = "'2.993.04"
$ws.Range("E16").Value = "  +4.89%  "

# Row 17
$ws.Range("D17").Value = "'25.99"
$ws.Range("E17").Value = "  +13.74%  "

# Row 18
$ws.Range("D18").Value = "'2.550.61"
$ws.Range("E18").Value = "  +5.84%  "

# Row 19
$ws.Range("D19").Value = "'11.62"
$ws.Range("E19").Value = "  +11.90%  "

# Row 20
$ws.Range("D20").Value = "'7.89"
$ws.Range("E20").Value = "  +15.07%  "

# Row 21
$ws.Range("D21").Value = "'365.98"
$ws.Range("E21").Value = "  +10.70%  "

# Row 22
$ws.Range("D22").Value = "'2.22"
$ws.Range("E22").Value = "  +19.07%  "

# Row 23
$ws.Range("D23").Value = "'4.04"
$ws.Range("E23").Value = "  +6.63%  "

# Row 24
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").Value = "'69.27"
$ws.Range("E25").Value = "  +4.76%  "

# Row 26
$ws.Range("D26").Value = "'4.13"
$ws.Range("E26").Value = "  +13.04%  "

# Row 27
$ws.Range("D27").Value = "'9.10"
$ws.Range("E27").Value = "  +11.47%  "

# Row 28
$ws.Range("D28").Value = "'2.724.14"
$ws.Range("E28").Value = "  +7.58%  "

# Row 29
$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -0.56%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0940"
$ws.Range("E30").Value = "  +16.21%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.88"
$ws.Range("E31").Value = "  +11.36%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'502.70"
$ws.Range("E32").Value = "  +19.27%  "

# Row 33
$ws.Range("D33").Value = "'1.34"
$ws.Range("E33").Value = "  +17.98%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.70"
$ws.Range("E34").Value = "  +6.40%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  -0.39%  "

# Row 36
$ws.Range("D36").Value = "'159.98"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").Value = "'0.118"
$ws.Range("E37").Value = "  +12.09%  "

# Row 38
$ws.Range("D38").Value = "'19.13"
$ws.Range("E38").Value = "  +7.46%  "

# Row 39
$ws.Range("E39").Value = "  +1.78%  "

# Row 40
$ws.Range("E40").Value = "  +0.05%  "

# Row 41
$ws.Range("D41").Value = "'4.85"
$ws.Range("E41").Value = "  +13.01%  "

# Row 42
$ws.Range("D42").Value = "'1.66"
$ws.Range("E42").Value = "  +12.77%  "

# Row 43
$ws.Range("D43").Value = "'0.318"
$ws.Range("E43").Value = "  +7.84%  "

# Row 44
$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  +20.12%  "

# Row 45
$ws.Range("D45").Value = "'38.94"
$ws.Range("E45").Value = "  +4.31%  "

# Row 46
$ws.Range("D46").Value = "'1.15"
$ws.Range("E46").Value = "  +7.73%  "

# Row 47
$ws.Range("D47").Value = "'148.15"
$ws.Range("E47").Value = "  +12.34%  "

# Row 48
$ws.Range("D48").Value = "'3.56"
$ws.Range("E48").Value = "  +7.52%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0783"
$ws.Range("E49").Value = "  +10.15%  "

# Row 50
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'0.519"
$ws.Range("E50").Value = "  +8.50%  "

# Row 51
$ws.Range("D51").Value = "'0.584"
$ws.Range("E51").Value = "  +5.43%  "
